$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the placeholder "picture" text for the stub/no-photo rows (11-13),
# leaving the cells blank since there is no actual photo for them.
$ws.Range("E11").Value = $null
$ws.Range("E12").Value = $null
$ws.Range("E13").Value = $null

# Update the last active selection / scroll position recorded in the sheet view.
[void]$ws.Range("D16").Select()

# Update window size/position recorded for the workbook view (best-effort;
# this is just the last-saved window geometry on the author's machine).
$excel.ActiveWindow.Width = 23260
$excel.ActiveWindow.Height = 12580
$excel.ActiveWindow.Top = 500
$excel.ActiveWindow.Left = 0
